$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @(32, 12, 2000),
    @(32, 8, 1600),
    @(32, 10, 2000),
    @(43, 9, 1959.2),
    @(43, 8, 2162.125),
    @(43, 13, -1890.2),
    @(43, 11, 1959.2),
    @(87, 8, 145000),
    @(90, 8, 145000),
    @(107, 8, 854),
    @(107, 9, 1094.4),
    @(107, 13, 825.5999999999999),
    @(107, 11, 1094.4),
    @(112, 10, 2622.7886),
    @(112, 8, 2611.0188),
    @(112, 14, -10084.3658),
    @(112, 12, 7868.3658),
    @(116, 14, -12047.4614),
    @(116, 8, 5112.294),
    @(116, 13, -1504),
    @(116, 11, 4946),
    @(116, 9, 4946),
    @(116, 12, 5163.4614),
    @(116, 10, 5163.4614),
    @(132, 11, 6535.071599999999),
    @(132, 13, -4005.071599999999),
    @(132, 9, 2178.3572),
    @(132, 8, 2499.7834),
    @(137, 9, 2417.3333),
    @(137, 11, 7251.999899999999),
    @(137, 8, 2391.25),
    @(137, 13, -4701.999899999999),
    @(138, 10, 3707.2888),
    @(138, 12, 11121.8664),
    @(138, 13, -514.0769),
    @(138, 11, 5654.0769),
    @(138, 8, 3298.776),
    @(138, 14, -21401.8664),
    @(138, 9, 1884.6923),
    @(32, 14, -2652)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @(26, 9, 7081.2),
    @(26, 14, -505660),
    @(26, 10, 505000),
    @(26, 8, 149343.72),
    @(26, 12, 505000),
    @(26, 11, 7081.2),
    @(26, 13, -6751.2),
    @(32, 12, 18641.555),
    @(32, 9, 1851.7675),
    @(32, 14, -19215.555),
    @(32, 11, 1851.7675),
    @(32, 8, 4757.6924),
    @(32, 10, 18641.555),
    @(32, 13, -1564.7675),
    @(61, 11, 4030.2632),
    @(61, 8, 5154.885),
    @(61, 13, -3818.2632),
    @(61, 9, 4030.2632),
    @(102, 13, -14080.286),
    @(102, 11, 15702.286),
    @(102, 8, 21322.133),
    @(102, 9, 15702.286),
    @(132, 11, 4499.3145),
    @(132, 13, -1969.3145),
    @(132, 9, 1499.7715),
    @(132, 8, 1827.619),
    @(136, 8, 5154.885),
    @(136, 13, -9540.7896),
    @(136, 11, 12090.7896),
    @(136, 9, 4030.2632)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @(99, 9, 2730.8096),
    @(99, 11, 2730.8096),
    @(99, 8, 2730.8096),
    @(99, 13, -1232.8096),
    @(105, 8, 3879.3333),
    @(105, 9, 2987.125),
    @(105, 11, 2987.125),
    @(105, 13, -1240.125),
    @(107, 12, 1288.75),
    @(107, 8, 2106.4075),
    @(107, 10, 1288.75),
    @(107, 9, 2248.6086),
    @(107, 14, -5128.75),
    @(107, 13, -328.6086),
    @(107, 11, 2248.6086)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @(31, 13, -6786.5713),
    @(31, 11, 7081.5713),
    @(31, 9, 7081.5713),
    @(31, 8, 6127.364),
    @(34, 13, -6879.5713),
    @(34, 9, 7081.5713),
    @(34, 11, 7081.5713),
    @(34, 8, 6127.364),
    @(58, 11, 2523.5715),
    @(58, 13, -2320.5715),
    @(58, 9, 2523.5715),
    @(58, 8, 2805.0833),
    @(107, 8, 29413492),
    @(107, 9, 41668332),
    @(107, 13, -41666412),
    @(107, 11, 41668332),
    @(122, 8, 2968.8064),
    @(122, 9, 3038.261),
    @(122, 13, -6664.782999999999),
    @(122, 11, 9114.782999999999),
    @(132, 11, 0),
    @(132, 9, 0),
    @(132, 14, -12510.5),
    @(132, 10, 2483.5),
    @(132, 12, 7450.5),
    @(132, 8, 2483.5),
    @(134, 9, 666.4286),
    @(134, 11, 1999.2858),
    @(134, 8, 773.3333),
    @(134, 13, 535.7142000000001),
    @(136, 8, 2805.0833),
    @(136, 13, -5020.7145),
    @(136, 11, 7570.7145),
    @(136, 9, 2523.5715)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}
$clears = @(
    @(132, 13)
)
foreach ($e in $clears) {
    $ws.Cells.Item($e[0], $e[1]).ClearContents()
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @(18, 8, 209),
    @(18, 11, 441),
    @(18, 9, 147),
    @(18, 13, -272),
    @(86, 8, 1179),
    @(86, 11, 2384.0001),
    @(86, 13, -1198.0001),
    @(86, 9, 794.6667),
    @(89, 8, 1179),
    @(89, 9, 794.6667),
    @(89, 11, 7152.0003),
    @(89, 13, -1224.0003),
    @(134, 9, 869.8461),
    @(134, 11, 2609.5383),
    @(134, 8, 4437.9375),
    @(134, 13, 2460.4617)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @(70, 9, 0),
    @(70, 12, 0),
    @(70, 11, 0),
    @(70, 8, 0),
    @(70, 10, 0),
    @(73, 10, 0),
    @(73, 8, 0),
    @(73, 11, 0),
    @(73, 12, 0),
    @(73, 9, 0),
    @(80, 11, 5168.8184),
    @(80, 8, 6468.524),
    @(80, 13, -4170.8184),
    @(80, 10, 7898.2),
    @(80, 9, 5168.8184),
    @(80, 14, -9894.200000000001),
    @(80, 12, 7898.2),
    @(83, 12, 39491),
    @(83, 10, 7898.2),
    @(83, 8, 6468.524),
    @(83, 11, 25844.092),
    @(83, 14, -49475),
    @(83, 9, 5168.8184),
    @(83, 13, -20852.092),
    @(97, 8, 402.86957),
    @(97, 11, 375.72726),
    @(97, 9, 375.72726),
    @(97, 13, 120.27274),
    @(102, 13, 588.0968),
    @(102, 12, 3524.6667),
    @(102, 11, 1033.9032),
    @(102, 8, 1437.8108),
    @(102, 10, 3524.6667),
    @(102, 9, 1033.9032),
    @(102, 14, -6768.6667),
    @(132, 11, 12063.3105),
    @(132, 13, -9533.3105),
    @(132, 9, 4021.1035),
    @(132, 14, -16224.1535),
    @(132, 10, 3721.3845),
    @(132, 12, 11164.1535),
    @(132, 8, 3928.3333)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}
$clears = @(
    @(70, 13),
    @(70, 14),
    @(73, 14),
    @(73, 13)
)
foreach ($e in $clears) {
    $ws.Cells.Item($e[0], $e[1]).ClearContents()
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @(7, 11, 2412.2222),
    @(7, 9, 2412.2222),
    @(7, 8, 2917.5),
    @(7, 13, -2300.2222),
    @(22, 14, -1689),
    @(22, 10, 1099),
    @(22, 12, 1099),
    @(22, 8, 1160),
    @(27, 8, 1160),
    @(27, 14, -1313),
    @(27, 10, 1099),
    @(27, 12, 1099),
    @(40, 13, -4793.4136),
    @(40, 9, 4929.4136),
    @(40, 11, 4929.4136),
    @(40, 8, 5064.5),
    @(46, 12, 1200),
    @(46, 10, 1200),
    @(46, 8, 1200),
    @(46, 14, -1576),
    @(122, 8, 1299.5),
    @(122, 9, 1299.5),
    @(122, 11, 3898.5),
    @(126, 8, 2917.5),
    @(126, 11, 7236.6666),
    @(126, 9, 2412.2222),
    @(126, 13, -4766.6666),
    @(122, 13, -1448.5)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @(113, 11, 1345.94118),
    @(113, 13, 824.05882),
    @(113, 8, 1110.762),
    @(113, 9, 448.64706),
    @(122, 8, 2945.5557),
    @(122, 9, 2087.1538),
    @(122, 13, -3811.4614),
    @(122, 11, 6261.4614),
    @(126, 14, -21963.5),
    @(126, 10, 5674.5),
    @(126, 8, 4180.909),
    @(126, 12, 17023.5),
    @(126, 11, 9982.2855),
    @(126, 9, 3327.4285),
    @(126, 13, -7512.2855),
    @(132, 11, 9038.5452),
    @(132, 13, -6508.5452),
    @(132, 9, 3012.8484),
    @(132, 8, 5247.8687)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}
